$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# swap F:V between rows 7 and 8
$ws.Cells.Item(7,6).Value = "Napredak"
$ws.Cells.Item(8,6).Value = "Crvena zvezda"
$ws.Cells.Item(7,7).Value = 0
$ws.Cells.Item(8,7).Value = 5
$ws.Cells.Item(7,8).Value = "Zeleznicar Pancevo"
$ws.Cells.Item(8,8).Value = "Vojvodina"
$ws.Cells.Item(7,9).Value = 1
$ws.Cells.Item(8,9).Value = 0
$ws.Cells.Item(7,10).Value = 1.83
$ws.Cells.Item(8,10).Value = 1.24
$ws.Cells.Item(7,11).Value = "29/07/2023 14:42"
$ws.Cells.Item(8,11).Value = "29/07/2023 08:12"
$ws.Cells.Item(7,12).Value = 2.05
$ws.Cells.Item(8,12).Value = 1.05
$ws.Cells.Item(7,13).Value = "30/07/2023 18:46"
$ws.Cells.Item(8,13).Value = "30/07/2023 18:03"
$ws.Cells.Item(7,14).Value = 3.18
$ws.Cells.Item(8,14).Value = 5.36
$ws.Cells.Item(7,15).Value = "29/07/2023 14:42"
$ws.Cells.Item(8,15).Value = "29/07/2023 08:12"
$ws.Cells.Item(7,16).Value = 3.11
$ws.Cells.Item(8,16).Value = 13.81
$ws.Cells.Item(7,17).Value = "30/07/2023 18:46"
$ws.Cells.Item(8,17).Value = "30/07/2023 19:07"
$ws.Cells.Item(7,18).Value = 3.9
$ws.Cells.Item(8,18).Value = 8.210000000000001
$ws.Cells.Item(7,19).Value = "29/07/2023 14:42"
$ws.Cells.Item(8,19).Value = "29/07/2023 08:12"
$ws.Cells.Item(7,20).Value = 3.83
$ws.Cells.Item(8,20).Value = 46.56
$ws.Cells.Item(7,21).Value = "30/07/2023 18:46"
$ws.Cells.Item(8,21).Value = "30/07/2023 19:07"
$ws.Cells.Item(7,22).Value = "https://www.betexplorer.com/football/serbia/super-liga/napredak-zeleznicar-pancevo/6savvMTp/"
$ws.Cells.Item(8,22).Value = "https://www.betexplorer.com/football/serbia/super-liga/crvena-zvezda-vojvodina/SIG0YBRm/"

# swap F:V between rows 19 and 20
$ws.Cells.Item(19,6).Value = "TSC"
$ws.Cells.Item(20,6).Value = "Vozdovac"
$ws.Cells.Item(19,7).Value = 1
$ws.Cells.Item(20,7).Value = 1
$ws.Cells.Item(19,8).Value = "Radnicki 1923"
$ws.Cells.Item(20,8).Value = "Radnik"
$ws.Cells.Item(19,9).Value = 0
$ws.Cells.Item(20,9).Value = 1
$ws.Cells.Item(19,10).Value = 1.33
$ws.Cells.Item(20,10).Value = 3.21
$ws.Cells.Item(19,11).Value = "11/08/2023 07:12"
$ws.Cells.Item(20,11).Value = "11/08/2023 07:12"
$ws.Cells.Item(19,12).Value = 1.47
$ws.Cells.Item(20,12).Value = 2.12
$ws.Cells.Item(19,13).Value = "12/08/2023 16:58"
$ws.Cells.Item(20,13).Value = "12/08/2023 18:54"
$ws.Cells.Item(19,14).Value = 4.58
$ws.Cells.Item(20,14).Value = 3.12
$ws.Cells.Item(19,15).Value = "11/08/2023 07:12"
$ws.Cells.Item(20,15).Value = "11/08/2023 07:12"
$ws.Cells.Item(19,16).Value = 4.35
$ws.Cells.Item(20,16).Value = 3.16
$ws.Cells.Item(19,17).Value = "12/08/2023 18:50"
$ws.Cells.Item(20,17).Value = "12/08/2023 18:54"
$ws.Cells.Item(19,18).Value = 7
$ws.Cells.Item(20,18).Value = 2.11
$ws.Cells.Item(19,19).Value = "11/08/2023 07:12"
$ws.Cells.Item(20,19).Value = "11/08/2023 07:12"
$ws.Cells.Item(19,20).Value = 6.26
$ws.Cells.Item(20,20).Value = 3.53
$ws.Cells.Item(19,21).Value = "12/08/2023 18:50"
$ws.Cells.Item(20,21).Value = "12/08/2023 18:54"
$ws.Cells.Item(19,22).Value = "https://www.betexplorer.com/football/serbia/super-liga/tsc-backa-topola-radnicki-1923/vVDe3xcd/"
$ws.Cells.Item(20,22).Value = "https://www.betexplorer.com/football/serbia/super-liga/fk-vozdovac-radnik-surdulica/WOOL9vKS/"

# swap F:V between rows 31 and 32
$ws.Cells.Item(31,6).Value = "Javor"
$ws.Cells.Item(32,6).Value = "Napredak"
$ws.Cells.Item(31,7).Value = 0
$ws.Cells.Item(32,7).Value = 0
$ws.Cells.Item(31,8).Value = "TSC"
$ws.Cells.Item(32,8).Value = "Mladost"
$ws.Cells.Item(31,9).Value = 3
$ws.Cells.Item(32,9).Value = 1
$ws.Cells.Item(31,10).Value = 4.58
$ws.Cells.Item(32,10).Value = 1.65
$ws.Cells.Item(31,11).Value = "17/08/2023 09:13"
$ws.Cells.Item(32,11).Value = "17/08/2023 09:13"
$ws.Cells.Item(31,12).Value = 3.83
$ws.Cells.Item(32,12).Value = 2.06
$ws.Cells.Item(31,13).Value = "20/08/2023 19:29"
$ws.Cells.Item(32,13).Value = "20/08/2023 19:20"
$ws.Cells.Item(31,14).Value = 3.69
$ws.Cells.Item(32,14).Value = 3.67
$ws.Cells.Item(31,15).Value = "17/08/2023 09:13"
$ws.Cells.Item(32,15).Value = "17/08/2023 09:13"
$ws.Cells.Item(31,16).Value = 3.48
$ws.Cells.Item(32,16).Value = 3.04
$ws.Cells.Item(31,17).Value = "20/08/2023 19:29"
$ws.Cells.Item(32,17).Value = "20/08/2023 19:20"
$ws.Cells.Item(31,18).Value = 1.6
$ws.Cells.Item(32,18).Value = 4.2
$ws.Cells.Item(31,19).Value = "17/08/2023 09:13"
$ws.Cells.Item(32,19).Value = "17/08/2023 09:13"
$ws.Cells.Item(31,20).Value = 1.91
$ws.Cells.Item(32,20).Value = 3.88
$ws.Cells.Item(31,21).Value = "20/08/2023 19:29"
$ws.Cells.Item(32,21).Value = "20/08/2023 19:20"
$ws.Cells.Item(31,22).Value = "https://www.betexplorer.com/football/serbia/super-liga/javor-tsc-backa-topola/23ajPJJF/"
$ws.Cells.Item(32,22).Value = "https://www.betexplorer.com/football/serbia/super-liga/napredak-mladost-lucani/GGOZGeRe/"

# swap F:V between rows 43 and 44
$ws.Cells.Item(43,6).Value = "Javor"
$ws.Cells.Item(44,6).Value = "Crvena zvezda"
$ws.Cells.Item(43,7).Value = 1
$ws.Cells.Item(44,7).Value = 2
$ws.Cells.Item(43,8).Value = "Radnicki Nis"
$ws.Cells.Item(44,8).Value = "Novi Pazar"
$ws.Cells.Item(43,9).Value = 0
$ws.Cells.Item(44,9).Value = 1
$ws.Cells.Item(43,10).Value = 2.55
$ws.Cells.Item(44,10).Value = 1.05
$ws.Cells.Item(43,11).Value = "31/08/2023 09:13"
$ws.Cells.Item(44,11).Value = "31/08/2023 09:13"
$ws.Cells.Item(43,12).Value = 2.56
$ws.Cells.Item(44,12).Value = 1.02
$ws.Cells.Item(43,13).Value = "02/09/2023 09:45"
$ws.Cells.Item(44,13).Value = "02/09/2023 17:53"
$ws.Cells.Item(43,14).Value = 2.93
$ws.Cells.Item(44,14).Value = 11.81
$ws.Cells.Item(43,15).Value = "31/08/2023 09:13"
$ws.Cells.Item(44,15).Value = "31/08/2023 09:13"
$ws.Cells.Item(43,16).Value = 3.07
$ws.Cells.Item(44,16).Value = 17.5
$ws.Cells.Item(43,17).Value = "02/09/2023 19:05"
$ws.Cells.Item(44,17).Value = "02/09/2023 17:55"
$ws.Cells.Item(43,18).Value = 2.62
$ws.Cells.Item(44,18).Value = 18.16
$ws.Cells.Item(43,19).Value = "31/08/2023 09:13"
$ws.Cells.Item(44,19).Value = "31/08/2023 09:13"
$ws.Cells.Item(43,20).Value = 2.55
$ws.Cells.Item(44,20).Value = 46.82
$ws.Cells.Item(43,21).Value = "02/09/2023 09:45"
$ws.Cells.Item(44,21).Value = "02/09/2023 17:53"
$ws.Cells.Item(43,22).Value = "https://www.betexplorer.com/football/serbia/super-liga/javor-radnicki-nis/Of7GTYA1/"
$ws.Cells.Item(44,22).Value = "https://www.betexplorer.com/football/serbia/super-liga/crvena-zvezda-novi-pazar/AkCzZTns/"

# swap F:V between rows 90 and 91
$ws.Cells.Item(90,6).Value = "Radnicki Nis"
$ws.Cells.Item(91,6).Value = "IMT Novi Beograd"
$ws.Cells.Item(90,7).Value = 1
$ws.Cells.Item(91,7).Value = 1
$ws.Cells.Item(90,8).Value = "Sp. Subotica"
$ws.Cells.Item(91,8).Value = "Crvena zvezda"
$ws.Cells.Item(90,9).Value = 1
$ws.Cells.Item(91,9).Value = 2
$ws.Cells.Item(90,10).Value = 1.7
$ws.Cells.Item(91,10).Value = 8.15
$ws.Cells.Item(90,11).Value = "27/10/2023 06:42"
$ws.Cells.Item(91,11).Value = "27/10/2023 06:42"
$ws.Cells.Item(90,12).Value = 1.66
$ws.Cells.Item(91,12).Value = 24.2
$ws.Cells.Item(90,13).Value = "28/10/2023 18:23"
$ws.Cells.Item(91,13).Value = "28/10/2023 18:29"
$ws.Cells.Item(90,14).Value = 3.42
$ws.Cells.Item(91,14).Value = 5.6
$ws.Cells.Item(90,15).Value = "27/10/2023 06:42"
$ws.Cells.Item(91,15).Value = "27/10/2023 06:42"
$ws.Cells.Item(90,16).Value = 3.65
$ws.Cells.Item(91,16).Value = 9.529999999999999
$ws.Cells.Item(90,17).Value = "28/10/2023 18:23"
$ws.Cells.Item(91,17).Value = "28/10/2023 18:29"
$ws.Cells.Item(90,18).Value = 4.23
$ws.Cells.Item(91,18).Value = 1.23
$ws.Cells.Item(90,19).Value = "27/10/2023 06:42"
$ws.Cells.Item(91,19).Value = "27/10/2023 06:42"
$ws.Cells.Item(90,20).Value = 5.08
$ws.Cells.Item(91,20).Value = 1.09
$ws.Cells.Item(90,21).Value = "28/10/2023 18:23"
$ws.Cells.Item(91,21).Value = "28/10/2023 18:21"
$ws.Cells.Item(90,22).Value = "https://www.betexplorer.com/football/serbia/super-liga/radnicki-nis-spartak-subotica/2qDshl5f/"
$ws.Cells.Item(91,22).Value = "https://www.betexplorer.com/football/serbia/super-liga/imt-novi-beograd-crvena-zvezda/SjAgknkD/"
# Append new rows 113-118 (matches rows 112..117 in new numbering / index 112..117)
$ws.Range("A112:V112").Copy()
$ws.Range("A113:V113").PasteSpecial(-4122)
$ws.Cells.Item(113,1).Value = 112
$ws.Cells.Item(113,2).Value = "serbia"
$ws.Cells.Item(113,3).Value = "super-liga"
$ws.Cells.Item(113,4).Value = "2023-2024"
$ws.Cells.Item(113,5).Value = 45255.58333333334
$ws.Cells.Item(113,6).Value = "Radnik"
$ws.Cells.Item(113,7).Value = 1
$ws.Cells.Item(113,8).Value = "Mladost"
$ws.Cells.Item(113,9).Value = 1
$ws.Cells.Item(113,10).Value = 2.3
$ws.Cells.Item(113,11).Value = "24/11/2023 02:12"
$ws.Cells.Item(113,12).Value = 2.47
$ws.Cells.Item(113,13).Value = "25/11/2023 13:56"
$ws.Cells.Item(113,14).Value = 2.87
$ws.Cells.Item(113,15).Value = "24/11/2023 02:12"
$ws.Cells.Item(113,16).Value = 2.84
$ws.Cells.Item(113,17).Value = "25/11/2023 13:56"
$ws.Cells.Item(113,18).Value = 3.02
$ws.Cells.Item(113,19).Value = "24/11/2023 02:12"
$ws.Cells.Item(113,20).Value = 3.18
$ws.Cells.Item(113,21).Value = "25/11/2023 13:56"
$ws.Cells.Item(113,22).Value = "https://www.betexplorer.com/football/serbia/super-liga/radnik-surdulica-mladost-lucani/ziUAE2IO/"

$ws.Range("A113:V113").Copy()
$ws.Range("A114:V114").PasteSpecial(-4122)
$ws.Cells.Item(114,1).Value = 113
$ws.Cells.Item(114,2).Value = "serbia"
$ws.Cells.Item(114,3).Value = "super-liga"
$ws.Cells.Item(114,4).Value = "2023-2024"
$ws.Cells.Item(114,5).Value = 45255.66666666666
$ws.Cells.Item(114,6).Value = "Vojvodina"
$ws.Cells.Item(114,7).Value = 1
$ws.Cells.Item(114,8).Value = "Crvena zvezda"
$ws.Cells.Item(114,9).Value = 2
$ws.Cells.Item(114,10).Value = 5.73
$ws.Cells.Item(114,11).Value = "24/11/2023 02:12"
$ws.Cells.Item(114,12).Value = 7.48
$ws.Cells.Item(114,13).Value = "25/11/2023 15:56"
$ws.Cells.Item(114,14).Value = 4.36
$ws.Cells.Item(114,15).Value = "24/11/2023 02:12"
$ws.Cells.Item(114,16).Value = 5.08
$ws.Cells.Item(114,17).Value = "25/11/2023 15:56"
$ws.Cells.Item(114,18).Value = 1.4
$ws.Cells.Item(114,19).Value = "24/11/2023 02:12"
$ws.Cells.Item(114,20).Value = 1.35
$ws.Cells.Item(114,21).Value = "25/11/2023 15:56"
$ws.Cells.Item(114,22).Value = "https://www.betexplorer.com/football/serbia/super-liga/vojvodina-crvena-zvezda/SQQ2GOmC/"

$ws.Range("A114:V114").Copy()
$ws.Range("A115:V115").PasteSpecial(-4122)
$ws.Cells.Item(115,1).Value = 114
$ws.Cells.Item(115,2).Value = "serbia"
$ws.Cells.Item(115,3).Value = "super-liga"
$ws.Cells.Item(115,4).Value = "2023-2024"
$ws.Cells.Item(115,5).Value = 45255.77083333334
$ws.Cells.Item(115,6).Value = "Partizan"
$ws.Cells.Item(115,7).Value = 0
$ws.Cells.Item(115,8).Value = "TSC"
$ws.Cells.Item(115,9).Value = 4
$ws.Cells.Item(115,10).Value = 1.52
$ws.Cells.Item(115,11).Value = "24/11/2023 02:12"
$ws.Cells.Item(115,12).Value = 1.52
$ws.Cells.Item(115,13).Value = "25/11/2023 18:29"
$ws.Cells.Item(115,14).Value = 3.85
$ws.Cells.Item(115,15).Value = "24/11/2023 02:12"
$ws.Cells.Item(115,16).Value = 4.2
$ws.Cells.Item(115,17).Value = "25/11/2023 18:29"
$ws.Cells.Item(115,18).Value = 5.11
$ws.Cells.Item(115,19).Value = "24/11/2023 02:12"
$ws.Cells.Item(115,20).Value = 5.71
$ws.Cells.Item(115,21).Value = "25/11/2023 18:29"
$ws.Cells.Item(115,22).Value = "https://www.betexplorer.com/football/serbia/super-liga/partizan-tsc-backa-topola/UgmYV3ma/"

$ws.Range("A115:V115").Copy()
$ws.Range("A116:V116").PasteSpecial(-4122)
$ws.Cells.Item(116,1).Value = 115
$ws.Cells.Item(116,2).Value = "serbia"
$ws.Cells.Item(116,3).Value = "super-liga"
$ws.Cells.Item(116,4).Value = "2023-2024"
$ws.Cells.Item(116,5).Value = 45256.54166666666
$ws.Cells.Item(116,6).Value = "Radnicki 1923"
$ws.Cells.Item(116,7).Value = 0
$ws.Cells.Item(116,8).Value = "Novi Pazar"
$ws.Cells.Item(116,9).Value = 4
$ws.Cells.Item(116,10).Value = 2
$ws.Cells.Item(116,11).Value = "24/11/2023 02:12"
$ws.Cells.Item(116,12).Value = 2.02
$ws.Cells.Item(116,13).Value = "26/11/2023 12:31"
$ws.Cells.Item(116,14).Value = 3.16
$ws.Cells.Item(116,15).Value = "24/11/2023 02:12"
$ws.Cells.Item(116,16).Value = 3.16
$ws.Cells.Item(116,17).Value = "26/11/2023 12:31"
$ws.Cells.Item(116,18).Value = 3.35
$ws.Cells.Item(116,19).Value = "24/11/2023 02:12"
$ws.Cells.Item(116,20).Value = 3.85
$ws.Cells.Item(116,21).Value = "26/11/2023 12:31"
$ws.Cells.Item(116,22).Value = "https://www.betexplorer.com/football/serbia/super-liga/radnicki-1923-novi-pazar/8vaPXPIn/"

$ws.Range("A116:V116").Copy()
$ws.Range("A117:V117").PasteSpecial(-4122)
$ws.Cells.Item(117,1).Value = 116
$ws.Cells.Item(117,2).Value = "serbia"
$ws.Cells.Item(117,3).Value = "super-liga"
$ws.Cells.Item(117,4).Value = "2023-2024"
$ws.Cells.Item(117,5).Value = 45256.625
$ws.Cells.Item(117,6).Value = "Radnicki Nis"
$ws.Cells.Item(117,7).Value = 0
$ws.Cells.Item(117,8).Value = "Cukaricki"
$ws.Cells.Item(117,9).Value = 2
$ws.Cells.Item(117,10).Value = 2.8
$ws.Cells.Item(117,11).Value = "24/11/2023 02:12"
$ws.Cells.Item(117,12).Value = 2.62
$ws.Cells.Item(117,13).Value = "26/11/2023 14:59"
$ws.Cells.Item(117,14).Value = 2.97
$ws.Cells.Item(117,15).Value = "24/11/2023 02:12"
$ws.Cells.Item(117,16).Value = 3.19
$ws.Cells.Item(117,17).Value = "26/11/2023 14:59"
$ws.Cells.Item(117,18).Value = 2.38
$ws.Cells.Item(117,19).Value = "24/11/2023 02:12"
$ws.Cells.Item(117,20).Value = 2.66
$ws.Cells.Item(117,21).Value = "26/11/2023 14:59"
$ws.Cells.Item(117,22).Value = "https://www.betexplorer.com/football/serbia/super-liga/radnicki-nis-cukaricki/nVcTWqYh/"

$ws.Range("A117:V117").Copy()
$ws.Range("A118:V118").PasteSpecial(-4122)
$ws.Cells.Item(118,1).Value = 117
$ws.Cells.Item(118,2).Value = "serbia"
$ws.Cells.Item(118,3).Value = "super-liga"
$ws.Cells.Item(118,4).Value = "2023-2024"
$ws.Cells.Item(118,5).Value = 45256.77083333334
$ws.Cells.Item(118,6).Value = "Zeleznicar Pancevo"
$ws.Cells.Item(118,7).Value = 2
$ws.Cells.Item(118,8).Value = "Napredak"
$ws.Cells.Item(118,9).Value = 1
$ws.Cells.Item(118,10).Value = 2
$ws.Cells.Item(118,11).Value = "24/11/2023 02:12"
$ws.Cells.Item(118,12).Value = 2
$ws.Cells.Item(118,13).Value = "26/11/2023 18:29"
$ws.Cells.Item(118,14).Value = 3.13
$ws.Cells.Item(118,15).Value = "24/11/2023 02:12"
$ws.Cells.Item(118,16).Value = 3.37
$ws.Cells.Item(118,17).Value = "26/11/2023 18:29"
$ws.Cells.Item(118,18).Value = 3.38
$ws.Cells.Item(118,19).Value = "24/11/2023 02:12"
$ws.Cells.Item(118,20).Value = 3.63
$ws.Cells.Item(118,21).Value = "26/11/2023 18:29"
$ws.Cells.Item(118,22).Value = "https://www.betexplorer.com/football/serbia/super-liga/zeleznicar-pancevo-napredak/MmQ6Fr3I/"
